# Update "想去人数" (F column) figures on the 展览 (sheet1) and 全部类型 (sheet4)
# sheets to reflect the latest scraped numbers from gh-pages output run 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 12360
$ws1.Range("F8").Value = 62
$ws1.Range("F9").Value = 20
$ws1.Range("F12").Value = 12188
$ws1.Range("F14").Value = 4714
$ws1.Range("F15").Value = 134
$ws1.Range("F22").Value = 169

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 12360
$ws4.Range("F10").Value = 62
$ws4.Range("F11").Value = 20
$ws4.Range("F14").Value = 12188
$ws4.Range("F16").Value = 4714
$ws4.Range("F17").Value = 134
$ws4.Range("F24").Value = 169
